$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(88, 1).Value = "2024-10-29 00:00:00"
$ws.Cells.Item(88, 2).Value = 74050
$ws.Cells.Item(88, 3).Value = 10360.27
$ws.Cells.Item(88, 4).Value = 9168.379999999999
$ws.Cells.Item(88, 5).Value = 7.1372
